$wb = $excel.ActiveWorkbook

# --- Rename existing sheets (lowercase, and "Em" -> "tom") ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "emre"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "jan"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "tom"

# --- Fill in data for the "tom" sheet (was empty) ---
$ws3.Range("A1").Value = "groupComponent"
$ws3.Range("B1").Value = "nameComponent"
$ws3.Range("A2").Value = "Motherboard"
$ws3.Range("B2").Value = "X99 Rampage V Extreme"
$ws3.Range("A3").Value = "RAM"
$ws3.Range("B3").Value = "GEL316GB1600C9DC"

# --- Add a brand-new sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "fvjhtjhtjht"
$ws4.Range("A1").Value = "groupComponent"
$ws4.Range("B1").Value = "nameComponent"

# --- Make "tom" (3rd sheet) the active/selected tab ---
$ws3.Select()
